# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing header cell so the new headers match
# the formatting (bold, bordered, centered) of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is the same for every player row (2 through 41):
# 79 wins, 83 losses, 0 ties.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
